# Auto-generated script applying the Halicarnassus_Profits market-data refresh
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) across
# multiple rows on multiple worksheets, matching the scheduled runner diff.

$wb = $excel.ActiveWorkbook

# ---- Worksheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 4741.231
$ws.Range("I28").Value = 1843.5
$ws.Range("J28").Value = 7225
$ws.Range("K28").Value = 1843.5
$ws.Range("L28").Value = 7225
$ws.Range("M28").Value = -1358.5
$ws.Range("N28").Value = -8195
$ws.Range("H106").Value = 7309.8887
$ws.Range("I106").Value = 6972.875
$ws.Range("K106").Value = 6972.875
$ws.Range("M106").Value = -6341.875
$ws.Range("H112").Value = 1676.8889
$ws.Range("I112").Value = 2831.6667
$ws.Range("J112").Value = 1099.5
$ws.Range("K112").Value = 8495.000100000001
$ws.Range("L112").Value = 3298.5
$ws.Range("M112").Value = -7387.000100000001
$ws.Range("N112").Value = -5514.5
$ws.Range("H126").Value = 45000
$ws.Range("J126").Value = 45000
$ws.Range("L126").Value = 45000
$ws.Range("N126").Value = -54880
$ws.Range("H132").Value = 2964.1482
$ws.Range("I132").Value = 1066.7391
$ws.Range("J132").Value = 13874.25
$ws.Range("K132").Value = 3200.2173
$ws.Range("L132").Value = 41622.75
$ws.Range("M132").Value = -670.2173000000003
$ws.Range("N132").Value = -46682.75

# ---- Worksheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 965.5
$ws.Range("I2").Value = 989.1667
$ws.Range("J2").Value = 894.5
$ws.Range("K2").Value = 989.1667
$ws.Range("L2").Value = 894.5
$ws.Range("M2").Value = -876.1667
$ws.Range("N2").Value = -1120.5
$ws.Range("H32").Value = 504.44446
$ws.Range("I32").Value = 504.44446
$ws.Range("K32").Value = 504.44446
$ws.Range("M32").Value = -217.44446
$ws.Range("H44").Value = 13709.053
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 13709.053
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 13709.053
$ws.Range("N44").Value = -14685.053
$ws.Range("H55").Value = 58521.75
$ws.Range("I55").Value = 22600
$ws.Range("J55").Value = 94443.5
$ws.Range("K55").Value = 22600
$ws.Range("L55").Value = 94443.5
$ws.Range("M55").Value = -22285
$ws.Range("N55").Value = -95073.5
$ws.Range("H61").Value = 4495
$ws.Range("I61").Value = 4495
$ws.Range("K61").Value = 4495
$ws.Range("M61").Value = -4283
$ws.Range("H97").Value = 679.3
$ws.Range("I97").Value = 691.2857
$ws.Range("J97").Value = 651.3333
$ws.Range("K97").Value = 691.2857
$ws.Range("L97").Value = 651.3333
$ws.Range("M97").Value = -195.2857
$ws.Range("N97").Value = -1643.3333
$ws.Range("H102").Value = 3456.9524
$ws.Range("I102").Value = 1572.6666
$ws.Range("K102").Value = 1572.6666
$ws.Range("M102").Value = 49.33339999999998
$ws.Range("H110").Value = 682.5
$ws.Range("I110").Value = 682.5
$ws.Range("K110").Value = 682.5
$ws.Range("M110").Value = 1362.5
$ws.Range("H116").Value = 965.5
$ws.Range("I116").Value = 989.1667
$ws.Range("J116").Value = 894.5
$ws.Range("K116").Value = 989.1667
$ws.Range("L116").Value = 894.5
$ws.Range("M116").Value = 1304.8333
$ws.Range("N116").Value = -5482.5
$ws.Range("H136").Value = 4495
$ws.Range("I136").Value = 4495
$ws.Range("K136").Value = 13485
$ws.Range("M136").Value = -10935
$ws.Range("M44").ClearContents()

# ---- Worksheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 965.5
$ws.Range("I3").Value = 989.1667
$ws.Range("J3").Value = 894.5
$ws.Range("K3").Value = 989.1667
$ws.Range("L3").Value = 894.5
$ws.Range("M3").Value = -875.1667
$ws.Range("N3").Value = -1122.5
$ws.Range("H20").Value = 3500
$ws.Range("I20").Value = 3500
$ws.Range("K20").Value = 3500
$ws.Range("M20").Value = -3253
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# ---- Worksheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 73005.42999999999
$ws.Range("I16").Value = 92314.63
$ws.Range("J16").Value = 2205
$ws.Range("K16").Value = 92314.63
$ws.Range("L16").Value = 2205
$ws.Range("M16").Value = -92027.63
$ws.Range("N16").Value = -2779
$ws.Range("H113").Value = 73005.42999999999
$ws.Range("I113").Value = 92314.63
$ws.Range("J113").Value = 2205
$ws.Range("K113").Value = 92314.63
$ws.Range("L113").Value = 2205
$ws.Range("M113").Value = -90144.63
$ws.Range("N113").Value = -6545
$ws.Range("H132").Value = 2507.5
$ws.Range("I132").Value = 2507.5
$ws.Range("K132").Value = 7522.5
$ws.Range("M132").Value = -4992.5
$ws.Range("H134").Value = 1391.5714
$ws.Range("I134").Value = 1469
$ws.Range("K134").Value = 4407
$ws.Range("M134").Value = -1872

# ---- Worksheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 7039.4707
$ws.Range("I109").Value = 752.3333
$ws.Range("J109").Value = 14112.5
$ws.Range("K109").Value = 2256.9999
$ws.Range("L109").Value = 42337.5
$ws.Range("M109").Value = -1216.9999
$ws.Range("N109").Value = -44417.5
$ws.Range("H137").Value = 5167.5
$ws.Range("J137").Value = 5231.6665
$ws.Range("L137").Value = 15694.9995
$ws.Range("N137").Value = -25894.9995

# ---- Worksheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 18333.334
$ws.Range("J103").Value = 18333.334
$ws.Range("L103").Value = 18333.334
$ws.Range("N103").Value = -20677.334
$ws.Range("H132").Value = 1700
$ws.Range("I132").Value = 1333.3334
$ws.Range("K132").Value = 4000.0002
$ws.Range("M132").Value = -1470.0002

# ---- Worksheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2500.5
$ws.Range("I16").Value = 2500.5
$ws.Range("K16").Value = 2500.5
$ws.Range("M16").Value = -2330.5
$ws.Range("H40").Value = 4273.6665
$ws.Range("J40").Value = 5831.6665
$ws.Range("L40").Value = 5831.6665
$ws.Range("N40").Value = -6103.6665
$ws.Range("H93").Value = 1574.4286
$ws.Range("I93").Value = 1804.8
$ws.Range("J93").Value = 998.5
$ws.Range("K93").Value = 1804.8
$ws.Range("L93").Value = 998.5
$ws.Range("M93").Value = -556.8
$ws.Range("N93").Value = -3494.5
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("H136").Value = 4999.6665
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 4999.6665
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 14998.9995
$ws.Range("N136").Value = -20098.9995
$ws.Range("N119").ClearContents()
$ws.Range("M136").ClearContents()

# ---- Worksheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1375.4375
$ws.Range("I132").Value = 1204.7273
$ws.Range("K132").Value = 3614.1819
$ws.Range("M132").Value = -1084.1819
$ws.Range("H136").Value = 2761.7144
$ws.Range("J136").Value = 5582.5
$ws.Range("L136").Value = 16747.5
$ws.Range("N136").Value = -21847.5

